$wb = $excel.ActiveWorkbook

# Sheets "展览" (exhibitions) and "全部类型" (all types) both carry the same
# rows of event data; update the "想去人数" (want-to-go count) column F for
# the same five events on each sheet.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 695
    $ws.Range("F4").Value = 532
    $ws.Range("F9").Value = 4179
    $ws.Range("F10").Value = 4306
    $ws.Range("F11").Value = 10
}
